$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2" = 57
    "C3" = 40
    "C4" = 24
    "C5" = 99
    "C6" = 113
    "C7" = 111
    "C8" = 98
    "C9" = 25
    "C10" = 117
    "C11" = 153
    "C12" = 140
    "C13" = 137
    "C14" = 155
    "C15" = 125
    "C16" = 102
    "C17" = 159
    "C18" = 122
    "C19" = 85
    "C20" = 82
    "C21" = 152
    "C23" = 28
    "C24" = 70
    "C25" = 44
    "C26" = 30
    "C27" = 45
    "C28" = 23
    "C29" = 46
    "C30" = 27
    "C31" = 18
    "C32" = 1
    "C33" = 55
    "C34" = 72
    "C35" = 60
    "C36" = 61
    "C37" = 88
    "C38" = 39
    "C39" = 62
    "C40" = 74
    "C41" = 29
    "C42" = 86
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
